$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1): Wins / Losses / Ties in columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header style (bold font, thin box border, centered/top aligned)
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Season record values for every player row (2-48): 89 wins, 73 losses, 0 ties
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 89
    $ws.Cells.Item($r, 31).Value = 73
    $ws.Cells.Item($r, 32).Value = 0
}
